$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Corrections to decision tree results in column J ("Correct Class (from logic tree)")
$ws.Range("J5").Value  = "diagnose"
$ws.Range("J17").Value = "diagnose"
$ws.Range("J21").Value = "adjust/ calibrate"
$ws.Range("J23").Value = "adjust/ calibrate"
$ws.Range("J31").Value = "replace"
$ws.Range("J32").Value = "diagnose"
$ws.Range("J34").Value = "adjust/ calibrate"
$ws.Range("J36").Value = "adjust/ calibrate"
